# edit.ps1
# Applies the "Playtestable" spelling-split, the "Dex's" spelling-split,
# the Week-22 Playtest/polish trim (+ _GoBack bookmark relocation), and the
# "Easter" run-merge to the Project Milestones table.
#
# Strategy: Word's Range.InsertXML replaces the *whole paragraph* that
# contains the matched Range when the inserted fragment is block-level
# (wrapped in <w:p>). So for every edit we Find the distinguishing text,
# then InsertXML a complete replacement <w:p> (reproducing the original
# paragraph's w:p attributes exactly) for that paragraph.

function ConvertTo-WordPackageXml {
    param([string]$BodyFragment)

    return '<?xml version="1.0" standalone="yes"?>' +
        '<?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' + $BodyFragment + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData>' +
        '</pkg:part>' +
        '</pkg:package>'
}

function Set-ParagraphXml {
    param(
        [string]$FindText,
        [string]$ParagraphXml
    )

    $d = $word.ActiveDocument
    $r = $d.Content
    $r.Find.ClearFormatting()
    $found = $r.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $FindText"
    }
    $pkg = ConvertTo-WordPackageXml($ParagraphXml)
    $r.InsertXML($pkg)
}

# NOTE: this interpreter does not support binding PowerShell's named
# (-Param value) arguments, so every call below is positional:
#   Set-ParagraphXml <FindText> <ParagraphXml>

# ---------------------------------------------------------------------
# 1) "Platestable" -> "Playtestable" (spell-checked, split into runs)
# ---------------------------------------------------------------------
$p1 = '<w:p w14:paraId="61876FE1" w14:textId="7F503603" w:rsidR="003B37A8" w:rsidRPr="00F06B74" w:rsidRDefault="00892A2C" w:rsidP="00471742">' +
      '<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Pla</w:t></w:r>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>y</w:t></w:r>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>testable</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> Prototype with placeholder assets and a small non-procedural track</w:t></w:r>' +
      '</w:p>'
Set-ParagraphXml "Platestable Prototype with placeholder assets and a small non-procedural track" $p1

# ---------------------------------------------------------------------
# 2) "...such as Dex's voice..." -> split out "Dex's" with spell markers
# ---------------------------------------------------------------------
$p2 = '<w:p w14:paraId="5317EA90" w14:textId="5AA2A926" w:rsidR="00892A2C" w:rsidRPr="00F06B74" w:rsidRDefault="00A100F5" w:rsidP="00892A2C">' +
      '<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Work on visual</w:t></w:r>' +
      '<w:r w:rsidR="00892A2C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> assets + </w:t></w:r>' +
      '<w:r w:rsidR="00B53C62"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>build</w:t></w:r>' +
      '<w:r w:rsidR="00892A2C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
      '<w:r w:rsidR="00B53C62"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">sound library for other sound effects (such as </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Dex' + [char]0x2019 + 's</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> voice, collecting dust, etc.)</w:t></w:r>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> and make sure these work in the composer.</w:t></w:r>' +
      '</w:p>'
Set-ParagraphXml "sound library for other sound effects" $p2

# ---------------------------------------------------------------------
# 3) Week 22: trim "+ also look at porting to mobile" and the highlighted
#    "(help me Fraser...)" aside; the _GoBack bookmark relocates here.
# ---------------------------------------------------------------------
$p3 = '<w:p w14:paraId="5870A2A7" w14:textId="49F093BA" w:rsidR="00A100F5" w:rsidRPr="00F06B74" w:rsidRDefault="00A100F5" w:rsidP="00A100F5">' +
      '<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">Playtest and polish depending on feedback and where we' + [char]0x2019 + 're at </w:t></w:r>' +
      '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
      '<w:bookmarkEnd w:id="0"/>' +
      '</w:p>'
Set-ParagraphXml "also look at porting to mobile" $p3

# ---------------------------------------------------------------------
# 4) "E" + _GoBack bookmark + "aster" -> single run "Easter"
#    (the old bookmark location is removed since it now lives in #3)
# ---------------------------------------------------------------------
$p4 = '<w:p w14:paraId="43D9064C" w14:textId="54816895" w:rsidR="00A100F5" w:rsidRPr="00F06B74" w:rsidRDefault="001C0AA7" w:rsidP="00A100F5">' +
      '<w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/></w:rPr><w:t>Easter</w:t></w:r>' +
      '</w:p>'
Set-ParagraphXml "aster" $p4

Write-Host "All edits applied."
